$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three "ECs" sending-cluster rows (old rows 2-4); remaining rows shift up
$ws.Rows.Item(2).Delete() | Out-Null
$ws.Rows.Item(2).Delete() | Out-Null
$ws.Rows.Item(2).Delete() | Out-Null

# Refresh remaining rows (now rows 2-7) with the updated TPM-derived values
# Row 2
$ws.Range("A2").Value = 'FAPs'
$ws.Range("B2").Value = 'Hgf'
$ws.Range("C2").Value = 'Met'
$ws.Range("D2").Value = 'ECs'
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.583520999999999
$ws.Range("H2").Value = 25.750563
$ws.Range("I2").Value = 0.8910607110509009
$ws.Range("J2").Value = 0.8910607110509009
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.140567
$ws.Range("N2").Value = 0.421701
$ws.Range("O2").Value = 0.07810038533383065
$ws.Range("P2").Value = 0.07810038533383065
$ws.Range("Q2").Value = 1.206559796407
$ws.Range("R2").Value = 10.859038167663
$ws.Range("S2").Value = 0.06959218488891249
$ws.Range("T2").Value = 0.06959218488891249

# Row 3
$ws.Range("A3").Value = 'FAPs'
$ws.Range("B3").Value = 'Hgf'
$ws.Range("C3").Value = 'Met'
$ws.Range("D3").Value = 'FAPs'
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.583520999999999
$ws.Range("H3").Value = 25.750563
$ws.Range("I3").Value = 0.8910607110509009
$ws.Range("J3").Value = 0.8910607110509009
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05371366666666667
$ws.Range("N3").Value = 0.161141
$ws.Range("O3").Value = 0.02984383293631935
$ws.Range("P3").Value = 0.02984383293631935
$ws.Range("Q3").Value = 0.4610523858203333
$ws.Range("R3").Value = 4.149471472383
$ws.Range("S3").Value = 0.02659266699672102
$ws.Range("T3").Value = 0.02659266699672102

# Row 4
$ws.Range("A4").Value = 'FAPs'
$ws.Range("B4").Value = 'Hgf'
$ws.Range("C4").Value = 'Met'
$ws.Range("D4").Value = 'MuSCs'
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.583520999999999
$ws.Range("H4").Value = 25.750563
$ws.Range("I4").Value = 0.8910607110509009
$ws.Range("J4").Value = 0.8910607110509009
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.605544
$ws.Range("N4").Value = 4.816632
$ws.Range("O4").Value = 0.8920557817298499
$ws.Range("P4").Value = 0.8920557817298499
$ws.Range("Q4").Value = 13.781220640424
$ws.Range("R4").Value = 124.030985763816
$ws.Range("S4").Value = 0.7948758591652673
$ws.Range("T4").Value = 0.7948758591652673

# Row 5
$ws.Range("A5").Value = 'MuSCs'
$ws.Range("B5").Value = 'Hgf'
$ws.Range("C5").Value = 'Met'
$ws.Range("D5").Value = 'ECs'
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.049404
$ws.Range("H5").Value = 3.148212
$ws.Range("I5").Value = 0.1089392889490991
$ws.Range("J5").Value = 0.1089392889490991
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.140567
$ws.Range("N5").Value = 0.421701
$ws.Range("O5").Value = 0.07810038533383065
$ws.Range("P5").Value = 0.07810038533383065
$ws.Range("Q5").Value = 0.147511572068
$ws.Range("R5").Value = 1.327604148612
$ws.Range("S5").Value = 0.008508200444918155
$ws.Range("T5").Value = 0.008508200444918155

# Row 6
$ws.Range("A6").Value = 'MuSCs'
$ws.Range("B6").Value = 'Hgf'
$ws.Range("C6").Value = 'Met'
$ws.Range("D6").Value = 'FAPs'
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.049404
$ws.Range("H6").Value = 3.148212
$ws.Range("I6").Value = 0.1089392889490991
$ws.Range("J6").Value = 0.1089392889490991
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05371366666666667
$ws.Range("N6").Value = 0.161141
$ws.Range("O6").Value = 0.02984383293631935
$ws.Range("P6").Value = 0.02984383293631935
$ws.Range("Q6").Value = 0.05636733665466667
$ws.Range("R6").Value = 0.507306029892
$ws.Range("S6").Value = 0.003251165939598333
$ws.Range("T6").Value = 0.003251165939598333

# Row 7
$ws.Range("A7").Value = 'MuSCs'
$ws.Range("B7").Value = 'Hgf'
$ws.Range("C7").Value = 'Met'
$ws.Range("D7").Value = 'MuSCs'
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.049404
$ws.Range("H7").Value = 3.148212
$ws.Range("I7").Value = 0.1089392889490991
$ws.Range("J7").Value = 0.1089392889490991
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.605544
$ws.Range("N7").Value = 4.816632
$ws.Range("O7").Value = 0.8920557817298499
$ws.Range("P7").Value = 0.8920557817298499
$ws.Range("Q7").Value = 1.684864295776
$ws.Range("R7").Value = 15.163778661984
$ws.Range("S7").Value = 0.09717992256458256
$ws.Range("T7").Value = 0.09717992256458256

